$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.489.30"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.99%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.667.85"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.04%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9972"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.36%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.10"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.94%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9978"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.38%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3950"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +1.30%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3919"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.82%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "52.10"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +5.34%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.408"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9973"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.39%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08585"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.98%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "24.53"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.42%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.333"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.12%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001340"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +3.53%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.857"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +4.12%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.660.01"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.53%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "95.65"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.71%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06963"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.72%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "20.60"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.78%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.011"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.38%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9971"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.33%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.473.44"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.91%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.436"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.90%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.043"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +9.95%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.26%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "157.94"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.36%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "143.03"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.352"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.140"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -4.87%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.519"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +3.52%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.845.15"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.79%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.079"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +8.20%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08256"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.25%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "11.29"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +12.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02990"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.69%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.852"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.67%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2760"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.23%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.09275"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.64%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.7771"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.14%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "13.86"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +4.82%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.447"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.97%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.47"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.84%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.7140"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +3.09%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.69%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.144"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.90%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9970"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.40%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08468"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.11%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "136.41"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.05%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.452"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +12.37%  "
